$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "V4" sheet as a copy of "V3", placed at the very front.
# ---------------------------------------------------------------------------
$v3 = $wb.Worksheets.Item("V3")
$v3.Copy($wb.Worksheets.Item(1))
$v4 = $wb.Worksheets.Item(1)
$v4.Name = "V4"

# ---------------------------------------------------------------------------
# 2. Insert two new rows (checkpoints) after row 9.
# ---------------------------------------------------------------------------
$v4.Range("A10:A11").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 3. Update the header row.
# ---------------------------------------------------------------------------
$v4.Range("B1").Value = "V4"
$v4.Range("C1").Value = "terrot"
$v4.Range("E1").Value = "Best"

# ---------------------------------------------------------------------------
# 4. Row level value updates reflecting the new run's progress.
# ---------------------------------------------------------------------------
$v4.Range("B3").ClearContents()

$v4.Range("B4").Value = 476

$v4.Range("B6").Value = 2579
$v4.Range("C6").Value = 3038

$v4.Range("B9").Value = 2813
$v4.Range("C9").Value = 3285

# New checkpoint rows.
$v4.Range("A10").Value = "Checkpoint 759"
$v4.Range("B10").Value = 3099
$v4.Range("C10").Value = 3571

$v4.Range("A11").Value = "Checkpoint 936"
$v4.Range("B11").Value = 3157
$v4.Range("C11").Value = 3630

# The rows that used to be 10 and 11 (now 12 and 13 after the insert) no
# longer carry "achieved" figures for this WIP.
$v4.Range("B12").ClearContents()
$v4.Range("C12").ClearContents()
$v4.Range("C13").ClearContents()
$v4.Range("C14").ClearContents()

# ---------------------------------------------------------------------------
# 5. Fix up the view on the new sheet (selection moved to A12).
# ---------------------------------------------------------------------------
$v4.Range("A12").Select()

# ---------------------------------------------------------------------------
# 6. The old "V3" sheet is no longer the active tab; reset its view.
# ---------------------------------------------------------------------------
$v3.Select()
$v3.Range("A2").Select()
$v3.Range("E52").Select()
